# Penalty Reward System (unfinished) - partial edit
# 1) On "Weekly Quantity" sheet, remove the two weekly rows for
#    2024-03-10 (45361.99999999999, qty 60) and 2024-03-17 (45368.99999999999, qty 40),
#    shifting subsequent rows up.
# 2) On "Monthly Trend" sheet, update March 2024 requested quantity from 120 to 20.

$wb = $excel.ActiveWorkbook

$weekly = $wb.Worksheets.Item("Weekly Quantity")
$monthly = $wb.Worksheets.Item("Monthly Trend")

# Delete rows 6 and 7 (1-based, row 1 is the header) on the Weekly Quantity sheet.
# Deleting row 6 first shifts row 7 up into row 6, so deleting row 6 twice
# removes both original rows 6 and 7.
$weekly.Rows.Item(6).Delete()
$weekly.Rows.Item(6).Delete()

# Update the Monthly Trend sheet's requested quantity for March 2024 (row 5, col B)
$monthly.Cells.Item(5, 2).Value = 20
